# povestka.docx — "90% work is finished" pass:
#   - the header row plus every "section title" (gridSpan=5) row in the
#     table get an explicit trHeight of 0
#   - those same section-title rows get grey (808080) shading
#     (three of them didn't have it yet; one already had shading and now
#     additionally gets its cell vertically centered)

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$grey = 8421504  # 0x808080

# Rows touched by this pass: the column-header row (1) and the five
# section-title rows (2, 3, 5, 8, 9). The data rows (4, 6, 7, 10) are
# left untouched.
foreach ($i in 1, 2, 3, 5, 8, 9) {
    $tbl.Rows.Item($i).Height = 0
}

# Row 2: "ТАЙИНЛАНМОҚДА" — already grey-shaded; now also vertically centered.
$tbl.Rows.Item(2).Cells.Item(1).VerticalAlignment = 1

# Row 3: "ВАЗИРЛИК  МАРКАЗИЙ  АППАРАТИ" (1st occurrence) — add grey shading.
$tbl.Rows.Item(3).Cells.Item(1).Shading.BackgroundPatternColor = $grey

# Row 5: "ҚОРАҚАЛПОҒИСТОН РЕСПУБЛИКАСИ ИИВ" — add grey shading.
$tbl.Rows.Item(5).Cells.Item(1).Shading.BackgroundPatternColor = $grey

# Row 8: "ҚОЛДИРИЛМОҚДА" — add grey shading.
$tbl.Rows.Item(8).Cells.Item(1).Shading.BackgroundPatternColor = $grey

# Row 9: "ВАЗИРЛИК  МАРКАЗИЙ  АППАРАТИ" (2nd occurrence) — add grey shading.
$tbl.Rows.Item(9).Cells.Item(1).Shading.BackgroundPatternColor = $grey
